$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying macro regenerated the worker/period rows in ascending
# (instead of descending) period order. Re-write the "Periodo Mora" column
# for the unchanged worker blocks (YESSICA and JORDIN) so the shared string
# table is rebuilt in the same ascending order as the refreshed report.
$periods16_30 = @("1607","1608","1609","1610","2009","2010","2011","2012","2101","2102","2103","2104","2105","2106","2107")
for ($i = 0; $i -lt $periods16_30.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods16_30[$i]
}

# Update the two summary value cells that changed with the refreshed data
$ws.Range("E11").Value = 480757
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 16

# The last surviving row of the "ANA MILENA" block (old row 35, period 2505)
# needs to take on the formatting that used to belong to the removed NIT
# totals row (old row 36) -- i.e. the "closing" bottom-border style used to
# cap the table.
$ws.Range("B36:J36").Copy()
$ws.Range("B35:J35").PasteSpecial(-4122)

# Remove HERLING DAVID ACEVEDO ACOSTA (2 periods) and the first two periods
# of ANA MILENA CASTELLON CAICEDO (2507, 2506) -- only the last period
# (2505) for ANA MILENA remains in the table.
$ws.Rows("31:34").Delete()

# Remove the old NIT/total summary row that used to close the table (it has
# shifted up to row 32 after the previous deletion).
$ws.Rows("32:32").Delete()
